# Regenerate s_vals data to filter save games.
# Updates columns B-E (and the dependent sum column G) for rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
    3 = @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569)
    4 = @(0.00006486019690155054, 0.004309184025731883, 0.7127328510149897, 6.48142807727062)
    5 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = $b + $c + $d + $e
}
